$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at sheet row 29 (shifts existing rows 29..108 down to 30..109)
$ws.Rows(29).Insert()

# Populate the new row 29 with the new weekly price entry
$ws.Range("A29").Value = 7
$ws.Range("B29").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C29").Value = "Ñuble"
$ws.Range("D29").Value = 44525
$ws.Range("D29").NumberFormat = $ws.Range("D30").NumberFormat
$ws.Range("E29").Value = 16
$ws.Range("F29").Value = 100112045
$ws.Range("G29").Value = "Zapallo"
$ws.Range("H29").Value = "Paine"
$ws.Range("I29").Value = "1a (guarda)"
$ws.Range("J29").Value = 120
$ws.Range("K29").Value = 220
$ws.Range("L29").Value = 250
$ws.Range("M29").Value = 235
$ws.Range("N29").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O29").Value = "Región de O'Higgins"
$ws.Range("P29").Value = 235
$ws.Range("Q29").Value = 1
$ws.Range("R29").Value = "Hortaliza"
